$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing cell values per diff (F and G columns for rows 608-698)
$ws.Range("F608").Value = 46644
$ws.Range("G608").Value = 2946
$ws.Range("F621").Value = 56501
$ws.Range("G630").Value = 2971
$ws.Range("F637").Value = 43774
$ws.Range("F639").Value = 40735
$ws.Range("F641").Value = 34409
$ws.Range("G641").Value = 1387
$ws.Range("F646").Value = 36019
$ws.Range("F650").Value = 38064
$ws.Range("F652").Value = 35208
$ws.Range("F653").Value = 34145
$ws.Range("F656").Value = 52559
$ws.Range("F658").Value = 27230
$ws.Range("F663").Value = 37262
$ws.Range("G663").Value = 1157
$ws.Range("F664").Value = 26421
$ws.Range("F666").Value = 23979
$ws.Range("G666").Value = 774
$ws.Range("F668").Value = 3379
$ws.Range("F670").Value = 52634
$ws.Range("F671").Value = 32615
$ws.Range("F672").Value = 29829
$ws.Range("G672").Value = 580
$ws.Range("F674").Value = 28713
$ws.Range("F677").Value = 56192
$ws.Range("G677").Value = 797
$ws.Range("F678").Value = 33833
$ws.Range("F679").Value = 29364
$ws.Range("F680").Value = 28470
$ws.Range("F681").Value = 26429
$ws.Range("G681").Value = 578
$ws.Range("F682").Value = 12585
$ws.Range("G682").Value = 413
$ws.Range("F683").Value = 24283
$ws.Range("F684").Value = 57095
$ws.Range("G684").Value = 1205
$ws.Range("F685").Value = 34471
$ws.Range("G685").Value = 1029
$ws.Range("F686").Value = 34426
$ws.Range("G686").Value = 1140
$ws.Range("F687").Value = 31426
$ws.Range("F688").Value = 32024
$ws.Range("G688").Value = 1347
$ws.Range("F689").Value = 15718
$ws.Range("F690").Value = 27694
$ws.Range("G690").Value = 1543
$ws.Range("F691").Value = 62071
$ws.Range("G691").Value = 2805
$ws.Range("F692").Value = 41516
$ws.Range("G692").Value = 2683
$ws.Range("F693").Value = 39361
$ws.Range("G693").Value = 2720
$ws.Range("F694").Value = 37344
$ws.Range("G694").Value = 2771
$ws.Range("F695").Value = 36495
$ws.Range("G695").Value = 3078
$ws.Range("F696").Value = 17560
$ws.Range("G696").Value = 2180
$ws.Range("F697").Value = 28226
$ws.Range("G697").Value = 2961
$ws.Range("F698").Value = 67671
$ws.Range("G698").Value = 5684

# Append new rows 699-704
$ws.Range("A699").Value = 44593
$ws.Range("A699").NumberFormat = "yyyy-mm-dd"
$ws.Range("B699").Value = 1043008
$ws.Range("C699").Value = 38506
$ws.Range("D699").Value = 20224
$ws.Range("E699").Value = 17877
$ws.Range("F699").Value = 42268
$ws.Range("G699").Value = 4223

$ws.Range("A700").Value = 44594
$ws.Range("A700").NumberFormat = "yyyy-mm-dd"
$ws.Range("B700").Value = 1062396
$ws.Range("C700").Value = 35766
$ws.Range("D700").Value = 19388
$ws.Range("E700").Value = 17896
$ws.Range("F700").Value = 41148
$ws.Range("G700").Value = 4058

$ws.Range("A701").Value = 44595
$ws.Range("A701").NumberFormat = "yyyy-mm-dd"
$ws.Range("B701").Value = 1081929
$ws.Range("C701").Value = 37175
$ws.Range("D701").Value = 19533
$ws.Range("E701").Value = 17921
$ws.Range("F701").Value = 36953
$ws.Range("G701").Value = 3337

$ws.Range("A702").Value = 44596
$ws.Range("A702").NumberFormat = "yyyy-mm-dd"
$ws.Range("B702").Value = 1101506
$ws.Range("C702").Value = 37141
$ws.Range("D702").Value = 19577
$ws.Range("E702").Value = 17938
$ws.Range("F702").Value = 28196
$ws.Range("G702").Value = 3084

$ws.Range("A703").Value = 44597
$ws.Range("A703").NumberFormat = "yyyy-mm-dd"
$ws.Range("B703").Value = 1117055
$ws.Range("C703").Value = 28424
$ws.Range("D703").Value = 15549
$ws.Range("E703").Value = 17958
$ws.Range("F703").Value = 11644
$ws.Range("G703").Value = 1862

$ws.Range("A704").Value = 44598
$ws.Range("A704").NumberFormat = "yyyy-mm-dd"
$ws.Range("B704").Value = 1127020
$ws.Range("C704").Value = 17497
$ws.Range("D704").Value = 9965
$ws.Range("E704").Value = 17973
$ws.Range("F704").Value = 14509
$ws.Range("G704").Value = 2248

